# [Kadastro App] Kayıt silindi: 1
#
# Deletes the record with "Kayıt No" = 1 (Tarsus / 18-UYG. /
# EMİNE ALANLI KIRCILI (K.Mühendisi)) from the master "Kayitlar" sheet
# and from its corresponding unit sheet "Tarsus".
#
# On "Kayitlar" the record currently lives in row 2, so deleting that row
# shifts the remaining record (Kayıt No 25 / Anamur / ...) up into row 2.
# On "Tarsus" the record is the only data row (row 2); deleting it leaves
# just the header row behind.

$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows(2).Delete()

$wsTarsus = $wb.Worksheets.Item("Tarsus")
$wsTarsus.Rows(2).Delete()
